$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; lift protection so cell values can be updated,
# then restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A38).
$disclaimer = $ws.Range("A38").Value2
$disclaimer = $disclaimer -replace "2021-05-26", "2021-05-27"
$ws.Range("A38").Value = $disclaimer

# Update Weight (D) and Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.03542824100887749
$ws.Range("E2").Value = -0.01302803000394781
$ws.Range("D3").Value = 0.02035137114781932
$ws.Range("E3").Value = -0.000779727095516547
$ws.Range("D4").Value = 0.01933476894233471
$ws.Range("E4").Value = -0.002399040383846529
$ws.Range("D5").Value = 0.03792575208053658
$ws.Range("E5").Value = 0.001046025104602499
$ws.Range("D6").Value = 0.03419008991370452
$ws.Range("E6").Value = 0
$ws.Range("D7").Value = 0.01976126072711013
$ws.Range("E7").Value = 0
$ws.Range("D8").Value = 0.03694948235198559
$ws.Range("E8").Value = -0.002986985278429666
$ws.Range("D9").Value = 0.02047770894441847
$ws.Range("E9").Value = -0.001877682403433556
$ws.Range("D10").Value = 0.02571258980301908
$ws.Range("E10").Value = -0.01087924043121358
$ws.Range("D11").Value = 0.02403164609222929
$ws.Range("E11").Value = 0.001851851851851771
$ws.Range("D12").Value = 0.05736102162113212
$ws.Range("E12").Value = 0.002604166666666741
$ws.Range("D13").Value = 0.02491967263354214
$ws.Range("E13").Value = -0.00257163850110198
$ws.Range("D14").Value = 0.02657488086724693
$ws.Range("E14").Value = -0.0009416195856873921
$ws.Range("D15").Value = 0.03210485508167212
$ws.Range("E15").Value = -0.001242015613910574
$ws.Range("D16").Value = 0.01920110721549793
$ws.Range("E16").Value = -0.00693730729701969
$ws.Range("D17").Value = 0.03198116203765879
$ws.Range("E17").Value = -0.009382951653943983
$ws.Range("D18").Value = 0.04196866329296263
$ws.Range("E18").Value = 0.001381533502187393
$ws.Range("D19").Value = 0.1255202628863726
$ws.Range("E19").Value = 0
$ws.Range("D20").Value = 0.009034322251252056
$ws.Range("E20").Value = -0.01696232034183609
$ws.Range("D21").Value = 0.01523392747615406
$ws.Range("E21").Value = -0.006800813293136088
$ws.Range("D22").Value = 0.01771368818911222
$ws.Range("E22").Value = -0.003904915312149093
$ws.Range("D23").Value = 0.01550048802042787
$ws.Range("E23").Value = -0.02288984263233185
$ws.Range("D24").Value = 0.02183507734845955
$ws.Range("E24").Value = -0.003405448717948678
$ws.Range("D25").Value = 0.01270376388218553
$ws.Range("E25").Value = -0.001341201716738349
$ws.Range("D26").Value = 0.04249939393205768
$ws.Range("E26").Value = -0.0069470623081479
$ws.Range("D27").Value = 0.02386614561310947
$ws.Range("E27").Value = 0.0000980296049406526
$ws.Range("D28").Value = 0.04556491648214644
$ws.Range("E28").Value = 0.003328578221588163
$ws.Range("D29").Value = 0.05615359034445588
$ws.Range("E29").Value = 0.003195455352387633
$ws.Range("D30").Value = 0.01325529651758052
$ws.Range("E30").Value = 0.01974522292993641
$ws.Range("D31").Value = 0.02054413292282365
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0.0135607145806151
$ws.Range("E32").Value = -0.01265822784810122
$ws.Range("D33").Value = 0.04165119126141194
$ws.Range("E33").Value = 0.001549586776859346
$ws.Range("D34").Value = 0.01708881453008761
$ws.Range("E34").Value = 0.003062563803412743
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = -0.001750165023660877

# Restore sheet protection.
$ws.Protect("")
